$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -6
    6  = -1
    7  = -1
    14 = -1
    15 = 1
    17 = 1
    18 = 0
    19 = -5
    20 = -1
    22 = 0
    25 = 0
    27 = -1
    35 = 1
    36 = 1
    41 = -1
    43 = -1
    44 = 2
    47 = 1
    48 = 1
    50 = 2
    51 = -10
    52 = -2
    54 = -3
    58 = 6
    60 = -3
    65 = -2
    66 = 3
    71 = -6
    74 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
